$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (row 5) down onto the new row 6,
# then set the new row's values explicitly.
$ws.Range("A5:F5").Copy() | Out-Null
$ws.Range("A6:F6").PasteSpecial(-4122) | Out-Null

$ws.Range("A6").Value = "EXL_CorporateLensHomePage_MyFavourites"
$ws.Range("B6").Value = "Add Favourites"
$ws.Range("C6").Value = "N"
$ws.Range("D6").Value = "Y"
$ws.Range("F6").Value = "Sprint1"

# Extend the Y/N list validation to include the new row.
$ws.Range("C2:D5").Validation.Delete() | Out-Null
$ws.Range("C2:D6").Validation.Add(3, 1, 1, '"Y,N"') | Out-Null

# Extend the Sprint list validation to include the new row.
$ws.Range("F2:F5").Validation.Delete() | Out-Null
$ws.Range("F2:F6").Validation.Add(3, 1, 1, '"Sprint1,Sprint2,Sprint3,Sprint4,Sprint5,Sprint6,Sprint7,Sprint8,Sprint9,Sprint10"') | Out-Null

# Move the active selection to the newly added row, as in the target workbook.
$ws.Range("A6").Select() | Out-Null
